# Apply the "2025 august output update" refresh to this FHIR StructureDefinition
# export workbook. The new canonical base URL (2rdoc.pt/ig/ios-lifestyle-medicine)
# replaces the old GitHub-shorthand based one everywhere it is used, the
# publication Date is bumped, and the bound ValueSet URL is refreshed too:
#   - Metadata!B2  (URL)               -> new canonical IG base URL for this StructureDefinition
#   - Metadata!B8  (Date)              -> new publication timestamp
#   - Elements!R5  (Fixed Value)       -> Extension.url is fixed to the SD's own canonical
#                                          URL, so it must be kept in sync with Metadata!B2
#   - Elements!Z6  (Binding Value Set) -> new canonical IG base URL for the bound ValueSet
#
# All of these cells hold plain text (shared strings), so we assign plain
# strings rather than typed/date values to avoid Excel re-typing them.

$wb = $excel.ActiveWorkbook

$newStructureDefinitionUrl = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/exposure-conditions"

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B2").Value = $newStructureDefinitionUrl
$metadata.Range("B8").Value = "2025-08-20T10:40:04+01:00"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("R5").Value = $newStructureDefinitionUrl
$elements.Range("Z6").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/ValueSet/exposure-conditions-vs"
